$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2260.4492
$ws.Range("I15").Value = 2260.4492
$ws.Range("K15").Value = 6781.3476
$ws.Range("M15").Value = -6612.3476

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1642.1428
$ws.Range("J46").Value = 1799
$ws.Range("L46").Value = 5397
$ws.Range("N46").Value = -5635

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 1642.1428
$ws.Range("J60").Value = 1799
$ws.Range("L60").Value = 5397
$ws.Range("N60").Value = -6365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1253.8667
$ws.Range("I92").Value = 684
$ws.Range("J92").Value = 1633.7778
$ws.Range("K92").Value = 684
$ws.Range("L92").Value = 1633.7778
$ws.Range("M92").Value = 564
$ws.Range("N92").Value = -4129.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 5921.2085
$ws.Range("I135").Value = 1664.375
$ws.Range("J135").Value = 14434.875
$ws.Range("K135").Value = 14979.375
$ws.Range("L135").Value = 129913.875
$ws.Range("M135").Value = -12444.375
$ws.Range("N135").Value = -134983.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8592.632
$ws.Range("I137").Value = 3806.9333
$ws.Range("J137").Value = 11713.739
$ws.Range("K137").Value = 11420.7999
$ws.Range("L137").Value = 35141.217
$ws.Range("M137").Value = -8870.7999
$ws.Range("N137").Value = -40241.217

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5558.4526
$ws.Range("I138").Value = 6299.2
$ws.Range("J138").Value = 5386.186
$ws.Range("K138").Value = 18897.6
$ws.Range("L138").Value = 16158.558
$ws.Range("M138").Value = -13757.6
$ws.Range("N138").Value = -26438.558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7823123
$ws.Range("I32").Value = 8779350
$ws.Range("K32").Value = 8779350
$ws.Range("M32").Value = -8779063

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 22778334
$ws.Range("I61").Value = 38466224
$ws.Range("K61").Value = 38466224
$ws.Range("M61").Value = -38466012

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 22743474
$ws.Range("I74").Value = 62500896
$ws.Range("J74").Value = 24948.143
$ws.Range("K74").Value = 62500896
$ws.Range("L74").Value = 24948.143
$ws.Range("M74").Value = -62500022
$ws.Range("N74").Value = -26696.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 22743474
$ws.Range("I77").Value = 62500896
$ws.Range("J77").Value = 24948.143
$ws.Range("K77").Value = 312504480
$ws.Range("L77").Value = 124740.715
$ws.Range("M77").Value = -312500112
$ws.Range("N77").Value = -133476.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 8420.823
$ws.Range("I102").Value = 8974.467000000001
$ws.Range("K102").Value = 8974.467000000001
$ws.Range("M102").Value = -7352.467000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 22778334
$ws.Range("I136").Value = 38466224
$ws.Range("K136").Value = 115398672
$ws.Range("M136").Value = -115396122

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1977.2727
$ws.Range("I86").Value = 1770
$ws.Range("K86").Value = 1770
$ws.Range("M86").Value = -647

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1977.2727
$ws.Range("I89").Value = 1770
$ws.Range("K89").Value = 8850
$ws.Range("M89").Value = -3234

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 28286.182
$ws.Range("J96").Value = 59138.75
$ws.Range("L96").Value = 59138.75
$ws.Range("N96").Value = -64630.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 30262.865
$ws.Range("I134").Value = 2732.2903
$ws.Range("J134").Value = 172504.17
$ws.Range("K134").Value = 8196.8709
$ws.Range("L134").Value = 517512.51
$ws.Range("M134").Value = -5661.8709
$ws.Range("N134").Value = -522582.51

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 366915.34
$ws.Range("I31").Value = 2402.3704
$ws.Range("J31").Value = 632911.3
$ws.Range("K31").Value = 2402.3704
$ws.Range("L31").Value = 632911.3
$ws.Range("M31").Value = -2107.3704
$ws.Range("N31").Value = -633501.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 366915.34
$ws.Range("I34").Value = 2402.3704
$ws.Range("J34").Value = 632911.3
$ws.Range("K34").Value = 2402.3704
$ws.Range("L34").Value = 632911.3
$ws.Range("M34").Value = -2200.3704
$ws.Range("N34").Value = -633315.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 1884.4445
$ws.Range("I49").Value = 691.2
$ws.Range("J49").Value = 3376
$ws.Range("K49").Value = 2073.6
$ws.Range("L49").Value = 10128
$ws.Range("M49").Value = -1917.6
$ws.Range("N49").Value = -10440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2117.5
$ws.Range("I68").Value = 2426.6667
$ws.Range("J68").Value = 2040.2084
$ws.Range("K68").Value = 7280.000100000001
$ws.Range("L68").Value = 6120.6252
$ws.Range("M68").Value = -6469.000100000001
$ws.Range("N68").Value = -7742.6252

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2117.5
$ws.Range("I71").Value = 2426.6667
$ws.Range("J71").Value = 2040.2084
$ws.Range("K71").Value = 21840.0003
$ws.Range("L71").Value = 18361.8756
$ws.Range("M71").Value = -17784.0003
$ws.Range("N71").Value = -26473.8756

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 4374.643
$ws.Range("J103").Value = 6261.222
$ws.Range("L103").Value = 18783.666
$ws.Range("N103").Value = -20541.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 777.37933
$ws.Range("I107").Value = 742.7857
$ws.Range("K107").Value = 2228.3571
$ws.Range("M107").Value = -308.3571000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 5945.8
$ws.Range("J112").Value = 6676.3335
$ws.Range("L112").Value = 20029.0005
$ws.Range("N112").Value = -22245.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 994.9487
$ws.Range("J113").Value = 946.96
$ws.Range("L113").Value = 2840.88
$ws.Range("N113").Value = -7180.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 440657.4
$ws.Range("I128").Value = 440657.4
$ws.Range("K128").Value = 1321972.2
$ws.Range("M128").Value = -1316992.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7972.3335
$ws.Range("I137").Value = 8457.25
$ws.Range("J137").Value = 7584.4
$ws.Range("K137").Value = 25371.75
$ws.Range("L137").Value = 22753.2
$ws.Range("M137").Value = -20271.75
$ws.Range("N137").Value = -32953.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4353.375
$ws.Range("I97").Value = 4261
$ws.Range("K97").Value = 4261
$ws.Range("M97").Value = -3765

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2406
$ws.Range("I102").Value = 1985.1111
$ws.Range("K102").Value = 1985.1111
$ws.Range("M102").Value = -363.1111000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3623
$ws.Range("I16").Value = 3696.353
$ws.Range("J16").Value = 2999.5
$ws.Range("K16").Value = 3696.353
$ws.Range("L16").Value = 2999.5
$ws.Range("M16").Value = -3526.353
$ws.Range("N16").Value = -3339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 32259708
$ws.Range("I93").Value = 38463244
$ws.Range("K93").Value = 38463244
$ws.Range("M93").Value = -38461996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17038.5
$ws.Range("J54").Value = 17038.5
$ws.Range("L54").Value = 17038.5
$ws.Range("N54").Value = -18078.5
